# Generate Report for Handoff
#
# The localization-status report was regenerated; the "80e49444-..." file's
# "Latest Handoff Datetime" on the zh-cn status sheet moved forward from
# 2016-03-31 05:11:08 to 2016-03-31 05:12:26.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("zh-cn")

# Row 4 on the zh-cn sheet is the 80e49444-9e42-4422-bfc9-ffd54ea8cdff entry;
# column E is "Latest Handoff Datetime". Keep it as literal text (it already
# carries a date display format via the cell style) rather than letting it be
# parsed into a date serial.
$ws.Range("E4").Value = "2016-03-31 05:12:26"
